$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 668, pushing the existing rows
# (old 668-682) down to 670-684.
$ws.Rows.Item(668).Insert()
$ws.Rows.Item(668).Insert()

# New row 668 - weekly update, Española / Primera, Provincia de Quillota
$ws.Range("A668").Value = 3
$ws.Range("B668").Value = "Femacal de La Calera"
$ws.Range("C668").Value = "Coquimbo"
$ws.Range("D668").Value = 45239
$ws.Range("E668").Value = 5
$ws.Range("F668").Value = 100112013
$ws.Range("G668").Value = "Alcachofa"
$ws.Range("H668").Value = "Española"
$ws.Range("I668").Value = "Primera"
$ws.Range("J668").Value = 12000
$ws.Range("K668").Value = 300
$ws.Range("L668").Value = 300
$ws.Range("M668").Value = 300
$ws.Range("N668").Value = "$/unidad"
$ws.Range("O668").Value = "Provincia de Quillota"
$ws.Range("P668").Value = 300
$ws.Range("Q668").Value = 1
$ws.Range("R668").Value = "Hortaliza"

# New row 669 - weekly update, Española / Segunda, Provincia de Quillota
$ws.Range("A669").Value = 3
$ws.Range("B669").Value = "Femacal de La Calera"
$ws.Range("C669").Value = "Coquimbo"
$ws.Range("D669").Value = 45239
$ws.Range("E669").Value = 5
$ws.Range("F669").Value = 100112013
$ws.Range("G669").Value = "Alcachofa"
$ws.Range("H669").Value = "Española"
$ws.Range("I669").Value = "Segunda"
$ws.Range("J669").Value = 12000
$ws.Range("K669").Value = 200
$ws.Range("L669").Value = 200
$ws.Range("M669").Value = 200
$ws.Range("N669").Value = "$/unidad"
$ws.Range("O669").Value = "Provincia de Quillota"
$ws.Range("P669").Value = 200
$ws.Range("Q669").Value = 1
$ws.Range("R669").Value = "Hortaliza"
